$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "246.05"
Set-TextValue $ws "D3" "22.20"
Set-TextValue $ws "D4" "5.359"
Set-TextValue $ws "D5" "0.05857"
Set-TextValue $ws "D7" "6.381"
Set-TextValue $ws "D9" "1.009"
Set-TextValue $ws "D10" "0.1424"
Set-TextValue $ws "D11" "0.03870"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws "D12" "0.07332"
Set-TextValue $ws "D13" "0.03004"
Set-TextValue $ws "D14" "4.179"
Set-TextValue $ws "D15" "0.09399"
Set-TextValue $ws "D16" "0.001585"
Set-TextValue $ws "D17" "0.04818"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D18" "0.005973"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D19" "0.004082"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D20" "0.0009877"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D21" "0.0001410"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D22" "3.689"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D23" "2.232"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D24" "0.01119"
$ws.Range("E24").Value = "23OneONEBestin24h"
Set-TextValue $ws "D26" "0.1296"
Set-TextValue $ws "D27" "0.0002472"
Set-TextValue $ws "D40" "0.03864"
Set-TextValue $ws "D41" "0.006371"
Set-TextValue $ws "D42" "0.1071"
Set-TextValue $ws "D43" "0.003000"
Set-TextValue $ws "D44" "0.005194"
Set-TextValue $ws "D45" "0.00005653"
Set-TextValue $ws "D47" "0.7221"
Set-TextValue $ws "D48" "0.07078"
